$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The publication date in C2 was stored as plain quoted text ("2020-03-11").
# Convert it to a real date value, formatted as yyyy-mm-dd, so it is a
# proper number instead of a shared-string.
$dateCell = $ws.Range("C2")
$dateCell.Style = "Normal"
$dateCell.Value2 = 43901
$dateCell.NumberFormat = "yyyy\-mm\-dd;@"

# The active selection moved from A2 to C2.
$dateCell.Select() | Out-Null
